$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H17").Value = 1362.5555
$ws.Range("J17").Value = 1369.4706
$ws.Range("L17").Value = 4108.4118
$ws.Range("N17").Value = -4444.4118
$ws.Range("H39").Value = 204.33333
$ws.Range("I39").Value = 218.45454
$ws.Range("J39").Value = 49
$ws.Range("K39").Value = 655.3636200000001
$ws.Range("L39").Value = 147
$ws.Range("M39").Value = -359.3636200000001
$ws.Range("N39").Value = -739
$ws.Range("H74").Value = 16666.6
$ws.Range("I74").Value = 17142.857
$ws.Range("K74").Value = 17142.857
$ws.Range("M74").Value = -16206.857
$ws.Range("H77").Value = 16666.6
$ws.Range("I77").Value = 17142.857
$ws.Range("K77").Value = 85714.285
$ws.Range("M77").Value = -81034.285
$ws.Range("H98").Value = 2774.577
$ws.Range("J98").Value = 2000
$ws.Range("L98").Value = 2000
$ws.Range("N98").Value = -4996
$ws.Range("H122").Value = 2774.577
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900
$ws.Range("H135").Value = 420.33334
$ws.Range("I135").Value = 444.6316
$ws.Range("K135").Value = 4001.6844
$ws.Range("M135").Value = -1466.6844
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1626.3636
$ws.Range("I32").Value = 1626.3636
$ws.Range("K32").Value = 1626.3636
$ws.Range("M32").Value = -1339.3636
$ws.Range("H62").Value = 72500
$ws.Range("J62").Value = 72500
$ws.Range("L62").Value = 72500
$ws.Range("N62").Value = -73748
$ws.Range("H65").Value = 72500
$ws.Range("J65").Value = 72500
$ws.Range("L65").Value = 217500
$ws.Range("N65").Value = -223740
$ws.Range("H74").Value = 2591.75
$ws.Range("I74").Value = 2117.5134
$ws.Range("K74").Value = 2117.5134
$ws.Range("M74").Value = -1243.5134
$ws.Range("H77").Value = 2591.75
$ws.Range("I77").Value = 2117.5134
$ws.Range("K77").Value = 10587.567
$ws.Range("M77").Value = -6219.566999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 31258134
$ws.Range("I20").Value = 38470910
$ws.Range("J20").Value = 2767
$ws.Range("K20").Value = 38470910
$ws.Range("L20").Value = 2767
$ws.Range("M20").Value = -38470663
$ws.Range("N20").Value = -3261
$ws.Range("H64").Value = 1550.75
$ws.Range("I64").Value = 1006
$ws.Range("J64").Value = 1732.3334
$ws.Range("K64").Value = 1006
$ws.Range("L64").Value = 1732.3334
$ws.Range("M64").Value = -781
$ws.Range("N64").Value = -2182.3334
$ws.Range("H67").Value = 1550.75
$ws.Range("I67").Value = 1006
$ws.Range("J67").Value = 1732.3334
$ws.Range("K67").Value = 1006
$ws.Range("L67").Value = 1732.3334
$ws.Range("M67").Value = -226
$ws.Range("N67").Value = -3292.3334
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 19999
$ws.Range("I56").Value = 19999
$ws.Range("K56").Value = 19999
$ws.Range("M56").Value = -19154
$ws.Range("H132").Value = 13339772
$ws.Range("I132").Value = 3409.5293
$ws.Range("K132").Value = 10228.5879
$ws.Range("M132").Value = -7698.5879
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 759.1429000000001
$ws.Range("I12").Value = 688
$ws.Range("J12").Value = 771
$ws.Range("K12").Value = 2064
$ws.Range("L12").Value = 2313
$ws.Range("M12").Value = -1891
$ws.Range("N12").Value = -2659
$ws.Range("H51").Value = 321.83334
$ws.Range("I51").Value = 95.333336
$ws.Range("J51").Value = 548.3333
$ws.Range("K51").Value = 286.000008
$ws.Range("L51").Value = 1644.9999
$ws.Range("M51").Value = 173.999992
$ws.Range("N51").Value = -2564.9999
$ws.Range("H56").Value = 7133.846
$ws.Range("I56").Value = 7133.846
$ws.Range("K56").Value = 7133.846
$ws.Range("M56").Value = -6603.846
$ws.Range("H68").Value = 12509753
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 12509753
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 37529259
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -37530881
$ws.Range("H71").Value = 12509753
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 12509753
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 112587777
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -112595889
$ws.Range("H75").Value = 64.5
$ws.Range("I75").Value = 49
$ws.Range("K75").Value = 147
$ws.Range("M75").Value = 851
$ws.Range("H78").Value = 64.5
$ws.Range("I78").Value = 49
$ws.Range("K78").Value = 441
$ws.Range("M78").Value = 4551
$ws.Range("H92").Value = 1049.75
$ws.Range("J92").Value = 1049.75
$ws.Range("L92").Value = 3149.25
$ws.Range("N92").Value = -5645.25
$ws.Range("H94").Value = 2809.8
$ws.Range("J94").Value = 2200
$ws.Range("L94").Value = 6600
$ws.Range("N94").Value = -7952
$ws.Range("H114").Value = 3671.25
$ws.Range("I114").Value = 2174
$ws.Range("K114").Value = 6522
$ws.Range("M114").Value = -3268
$ws.Range("H132").Value = 1833.7693
$ws.Range("I132").Value = 1568.75
$ws.Range("K132").Value = 14118.75
$ws.Range("M132").Value = -11588.75
$ws.Range("H137").Value = 2662.7334
$ws.Range("J137").Value = 2989.7144
$ws.Range("L137").Value = 8969.143199999999
$ws.Range("N137").Value = -19169.1432
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7535.6
$ws.Range("I70").Value = 4999.3335
$ws.Range("K70").Value = 4999.3335
$ws.Range("M70").Value = -4729.3335
$ws.Range("H73").Value = 7535.6
$ws.Range("I73").Value = 4999.3335
$ws.Range("K73").Value = 4999.3335
$ws.Range("M73").Value = -4063.3335
$ws.Range("H80").Value = 142862290
$ws.Range("I80").Value = 250004500
$ws.Range("J80").Value = 5998
$ws.Range("K80").Value = 250004500
$ws.Range("L80").Value = 5998
$ws.Range("M80").Value = -250003502
$ws.Range("N80").Value = -7994
$ws.Range("H82").Value = 60000
$ws.Range("I82").Value = 50000
$ws.Range("J82").Value = 70000
$ws.Range("K82").Value = 50000
$ws.Range("L82").Value = 70000
$ws.Range("M82").Value = -49617
$ws.Range("N82").Value = -70766
$ws.Range("H83").Value = 142862290
$ws.Range("I83").Value = 250004500
$ws.Range("J83").Value = 5998
$ws.Range("K83").Value = 1250022500
$ws.Range("L83").Value = 29990
$ws.Range("M83").Value = -1250017508
$ws.Range("N83").Value = -39974
$ws.Range("H85").Value = 60000
$ws.Range("I85").Value = 50000
$ws.Range("J85").Value = 70000
$ws.Range("K85").Value = 50000
$ws.Range("L85").Value = 70000
$ws.Range("M85").Value = -48674
$ws.Range("N85").Value = -72652
$ws.Range("H132").Value = 2094.6667
$ws.Range("I132").Value = 1767.5927
$ws.Range("K132").Value = 5302.7781
$ws.Range("M132").Value = -2772.7781
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1933.3334
$ws.Range("I10").Value = 1933.3334
$ws.Range("K10").Value = 1933.3334
$ws.Range("M10").Value = -1793.3334
$ws.Range("H46").Value = 3296.8147
$ws.Range("I46").Value = 2440.4375
$ws.Range("J46").Value = 4542.4546
$ws.Range("K46").Value = 2440.4375
$ws.Range("L46").Value = 4542.4546
$ws.Range("M46").Value = -2252.4375
$ws.Range("N46").Value = -4918.4546
$ws.Range("H56").Value = 8700.5
$ws.Range("I56").Value = 8700.5
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 8700.5
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -8009.5
$ws.Range("N56").ClearContents()
$ws.Range("H132").Value = 8413.519
$ws.Range("I132").Value = 5363.1763
$ws.Range("K132").Value = 16089.5289
$ws.Range("M132").Value = -13559.5289
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 8083
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 8083
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 8083
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -8311
$ws.Range("H41").Value = 18952.5
$ws.Range("J41").Value = 18936.666
$ws.Range("L41").Value = 18936.666
$ws.Range("N41").Value = -19716.666
$ws.Range("H97").Value = 20000
$ws.Range("I97").Value = 20000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 20000
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -19009
$ws.Range("N97").ClearContents()
